# Append a new data row (row 7) to the "Artfynd" worksheet, mirroring the
# structure of the existing rows (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $text)
    # Force the cell to be stored as text even though its content could be
    # parsed as a number or date, then restore the default (unstyled) look.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$row = 7

$ws.Cells.Item($row, 1).Value = 74683551              # A7 Id
$ws.Cells.Item($row, 2).Value = 44331                 # B7 Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value = "Ovaliderad"          # C7 Valideringsstatus
$ws.Cells.Item($row, 4).Value = "NT"                  # D7 Rödlistade
$ws.Cells.Item($row, 5).Value = 201164                # E7 TaxonId
$ws.Cells.Item($row, 6).Value = "Sexfläckig bastardsvärmare"   # F7 Artnamn
$ws.Cells.Item($row, 7).Value = "Zygaena filipendulae"         # G7 Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Linnaeus, 1758)"             # H7 Auktor

Set-TextCell $ws.Cells.Item($row, 9) "1"              # I7 Antal (kept as text)

$ws.Cells.Item($row, 16).Value = "Klinte, vägrenen, Klinte, Gtl" # P7 Lokalnamn
$ws.Cells.Item($row, 17).Value = 730441               # Q7 Ost
$ws.Cells.Item($row, 18).Value = 6368156              # R7 Nord
$ws.Cells.Item($row, 19).Value = 50                   # S7 Noggrannhet
$ws.Cells.Item($row, 20).Value = "Gotland"            # T7 Län
$ws.Cells.Item($row, 21).Value = "Gotland"            # U7 Kommun
$ws.Cells.Item($row, 22).Value = "Gotland"            # V7 Provins
$ws.Cells.Item($row, 23).Value = "Gammelgarn"         # W7 Församling

Set-TextCell $ws.Cells.Item($row, 25) "2018-07-11"    # Y7 Startdatum (text)
Set-TextCell $ws.Cells.Item($row, 27) "2018-07-11"    # AA7 Slutdatum (text)

$ws.Cells.Item($row, 30).Value = $false               # AD7 Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false               # AE7 Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false               # AG7 Ospontan

$ws.Cells.Item($row, 49).Value = "Ulf L Larsson"      # AW7 Rapportör
$ws.Cells.Item($row, 50).Value = "Ulf L Larsson"      # AX7 Observatörer
